# Update gh-pages to output generated at 456a3b4
# Sheet "展览" (1st sheet) and sheet "全部类型" (4th sheet) contain overlapping
# rows of event data that both need the same refreshed counts/text.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

$newVenue = "幸福路1号(筑梦集团·结婚产业园·B1幢) 费加罗宴会艺术中心(旗舰店)"

# --- Sheet 1: 展览 ---
$ws1.Cells.Item(3, 6).Value2 = 7530    # F3: 7522 -> 7530
$ws1.Cells.Item(5, 6).Value2 = 14      # F5: 13 -> 14
$ws1.Cells.Item(6, 4).Value2 = $newVenue   # D6: venue text update
$ws1.Cells.Item(6, 6).Value2 = 454     # F6: 451 -> 454
$ws1.Cells.Item(7, 6).Value2 = 4138    # F7: 4121 -> 4138
$ws1.Cells.Item(9, 6).Value2 = 578     # F9: 577 -> 578
$ws1.Cells.Item(12, 6).Value2 = 148    # F12: 147 -> 148

# --- Sheet 4: 全部类型 ---
$ws4.Cells.Item(4, 6).Value2 = 7530    # F4: 7522 -> 7530
$ws4.Cells.Item(7, 6).Value2 = 14      # F7: 13 -> 14
$ws4.Cells.Item(8, 4).Value2 = $newVenue   # D8: venue text update
$ws4.Cells.Item(8, 6).Value2 = 454     # F8: 451 -> 454
$ws4.Cells.Item(9, 6).Value2 = 4138    # F9: 4121 -> 4138
$ws4.Cells.Item(11, 6).Value2 = 578    # F11: 577 -> 578
$ws4.Cells.Item(15, 6).Value2 = 148    # F15: 147 -> 148
